# Actualización automática de grupos experimentales
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the experimental group labels for several rows.
$ws.Range("B2").Value  = "Sin SmartScore"
$ws.Range("B3").Value  = "Con SmartScore"
$ws.Range("B4").Value  = "Con SmartScore"
$ws.Range("B6").Value  = "Sin SmartScore"
$ws.Range("B7").Value  = "Con SmartScore"
$ws.Range("B8").Value  = "Sin SmartScore"
$ws.Range("B10").Value = "Sin SmartScore"
$ws.Range("B13").Value = "Con SmartScore"
$ws.Range("B14").Value = "Con SmartScore"

# Row 14 SmartScore values were stored as text; convert them to real numbers.
$ws.Range("I14").Value  = 0.589
$ws.Range("L14").Value  = 0.52
$ws.Range("O14").Value  = 0.494
$ws.Range("R14").Value  = 0.562
$ws.Range("U14").Value  = 0.553
$ws.Range("X14").Value  = 0.545
$ws.Range("AA14").Value = 0.721
$ws.Range("AD14").Value = 0.622
$ws.Range("AG14").Value = 0.61
